# Update the cryptocurrency price ("Price") and hourly volume change ("Volume(1h)") values
# for the rows that changed in this data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.477.00"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.818.04"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.04"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5077"
$ws.Range("E7").Value = "  -5.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3852"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08481"
$ws.Range("E9").Value = "  +9.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.91"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.429"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.03"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.492"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "1.817.31"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001139"
$ws.Range("E17").Value = "  +4.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.13"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06685"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.079"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "28.510.95"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.42"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.23"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.53"
$ws.Range("D28").Value = "2.025.15"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.396"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.94"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.091"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1077"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.738"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.687"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07391"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2228"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.216"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.755"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6326"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.24"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.190"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.51"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.752"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5907"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.69"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.987"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.192"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06985"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.07"
$ws.Range("E51").Value = "  -0.42%  "
